$d = $word.ActiveDocument

# 1. Add a period after "Enseñar los avances de Juan Carlos con la BB.DD"
$r = $d.Content
$found = $r.Find.Execute("Enseñar los avances de Juan Carlos con la BB.DD", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)
if ($found) {
    $r.Collapse(0)
    $r.InsertAfter(".")
}

# 2. Append four new sub-bullet (level 2) paragraphs at the end of the document,
#    after "Tratar el tema de: ¿Qué cuestiones le planteamos al tutor el lunes 24?"
$newItems = @(
    "¿Nos recomiendas usar Laravel?",
    "¿Cómo hay que presentar la documentacion del proyecto de cara a su evaluacion, y de que cosas tiene que constar?",
    "¿Cada cuanto vas a hacer las sesiones de control para ver cómo vamos?",
    "Sugerencias para preguntarle más cosas […]"
)

$isFirstNewItem = $true
foreach ($itemText in $newItems) {
    $lastPara = $d.Paragraphs.Last
    $insertRange = $lastPara.Range
    $insertRange.Collapse(0)
    $insertRange.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Last
    $newPara.Range.Text = $itemText
    if ($isFirstNewItem) {
        $newPara.Range.ListFormat.ListIndent()
        $isFirstNewItem = $false
    }
}
